$wb = $excel.ActiveWorkbook

# "Pour la prochaine fois" is the first / active sheet in the workbook.
$ws = $wb.Worksheets.Item("Pour la prochaine fois")

# Delete the top block of rows (old rows 2-6: the "hitbox" note pair at
# B3:C3 / B4:C4, and the lone styled blank cell at B5). Everything below
# shifts up by 5 rows.
$ws.Rows("2:6").Delete()

# The "github pipeline ? / github projet/action" note (now sitting at row
# 17) is dropped entirely...
$ws.Range("B17:C17").ClearContents()

# Expand a couple of existing notes with extra wording.
$ws.Range("C10").Value = "pour une prochiane co, sauvegarder dans le temps ? Donc pouvoir se co avec un pseduo déjà enregistrée, mais pas connecté"

# ...and replaced by a brand-new single-cell note one row higher, at B16.
$ws.Range("B16").Value = "faire en sorte que le site est accessible avec le www."

$ws.Range("B18").Value = "quand on tue quelqu’un, deplacer le perso au spawn, reglé ???"

# Two brand-new tasks appended at the bottom of the list.
$ws.Range("B19").Value = "Faire en sorte que on peut pas jouer avec un pseudo, seulement si il est dans la game"
$ws.Range("B20").Value = "pseduo taille max"

$ws.Range("B24").Select()
